$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.460.30'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.729.16'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4869'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2622'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06176'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.739.00'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07021'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.49'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.560'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6006'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.20'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.482.61'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007080'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.19%  '
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.956.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.463'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.608'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.186'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.28'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.407'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '106.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.714'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.962'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07969'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.690'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04527'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.76%  '
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.613'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6239'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9091'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('E39').Value = '  -5.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.409'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01483'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.436'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3863'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.679'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1158'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05366'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.30'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.699'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('E51').Value = '  -1.18%  '
